# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 17: bring MAYERLIN CABARCAS CARRASQUILLA record up (was previously on row 21)
$ws.Range("C17").Value = "33102296"
$ws.Range("D17").Value = "MAYERLIN CABARCAS CARRASQUILLA"
$ws.Range("E17").Value = "2102"
$ws.Range("F17").Value = 36341
$ws.Range("G17").Value = 1000000

# Row 18: JULIA CLARISA MORALES DE HORTA, period 2205
$ws.Range("C18").Value = "1128044601"
$ws.Range("D18").Value = "JULIA CLARISA MORALES DE HORTA"
$ws.Range("E18").Value = "2205"
$ws.Range("F18").Value = 40000
$ws.Range("G18").Value = 1000000

# Row 19: JULIA CLARISA MORALES DE HORTA, period 2206
$ws.Range("C19").Value = "1128044601"
$ws.Range("D19").Value = "JULIA CLARISA MORALES DE HORTA"
$ws.Range("E19").Value = "2206"
$ws.Range("F19").Value = 40000
$ws.Range("G19").Value = 1000000

# Row 20: JULIA CLARISA MORALES DE HORTA, period 2207
$ws.Range("C20").Value = "1128044601"
$ws.Range("D20").Value = "JULIA CLARISA MORALES DE HORTA"
$ws.Range("E20").Value = "2207"
$ws.Range("F20").Value = 40000
$ws.Range("G20").Value = 1000000

# Row 21: JULIA CLARISA MORALES DE HORTA, period 2208
$ws.Range("C21").Value = "1128044601"
$ws.Range("D21").Value = "JULIA CLARISA MORALES DE HORTA"
$ws.Range("E21").Value = "2208"
$ws.Range("F21").Value = 16000
$ws.Range("G21").Value = 1000000
